$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") '27.980.52'
$ws.Range("E2").Value = '  +0.10%  '

# Row 3
Set-TextValue $ws.Range("D3") '1.858.12'
$ws.Range("E3").Value = '  -0.47%  '

# Row 4
Set-TextValue $ws.Range("D4") '1.003'
$ws.Range("E4").Value = '  +0.13%  '

# Row 5
Set-TextValue $ws.Range("D5") '311.46'
$ws.Range("E5").Value = '  -0.28%  '

# Row 6
Set-TextValue $ws.Range("D6") '1.002'
$ws.Range("E6").Value = '  +0.05%  '

# Row 7
Set-TextValue $ws.Range("D7") '0.5084'
$ws.Range("E7").Value = '  +2.02%  '

# Row 8
Set-TextValue $ws.Range("D8") '0.3809'
$ws.Range("E8").Value = '  -0.34%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.08262'
$ws.Range("E9").Value = '  -7.47%  '

# Row 10
Set-TextValue $ws.Range("D10") '1.109'
$ws.Range("E10").Value = '  -0.89%  '

# Row 11
$ws.Range("B11").Value = 'OKB'
$ws.Range("C11").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range("D11") '41.50'
$ws.Range("E11").Value = '  +0.01%  '

# Row 12
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range("D12") '6.200'
$ws.Range("E12").Value = '  -2.75%  '

# Row 13
$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue $ws.Range("D13") '20.50'
$ws.Range("E13").Value = '  -0.89%  '

# Row 14
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range("D14") '1.861.47'
$ws.Range("E14").Value = '  -0.15%  '

# Row 15
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range("D15") '7.201'
$ws.Range("E15").Value = '  -0.32%  '

# Row 16
$ws.Range("B16").Value = 'BinanceUSD'
$ws.Range("C16").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue $ws.Range("D16") '1.003'
$ws.Range("E16").Value = '  +0.09%  '

# Row 17
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Range("D17") '0.00001096'
$ws.Range("E17").Value = '  -0.22%  '

# Row 18
$ws.Range("B18").Value = 'Litecoin'
$ws.Range("C18").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue $ws.Range("D18") '90.51'
$ws.Range("E18").Value = '  -0.58%  '

# Row 19
$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws.Range("D19") '0.06597'
$ws.Range("E19").Value = '  -0.98%  '

# Row 20
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue $ws.Range("D20") '17.68'
$ws.Range("E20").Value = '  -2.14%  '

# Row 21
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws.Range("D21") '1.002'
$ws.Range("E21").Value = '  +0.15%  '

# Row 22
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws.Range("D22") '6.022'
$ws.Range("E22").Value = '  -1.60%  '

# Row 23
$ws.Range("B23").Value = 'WrappedBTC'
$ws.Range("C23").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue $ws.Range("D23") '27.992.62'
$ws.Range("E23").Value = '  +0.05%  '

# Row 24
$ws.Range("B24").Value = 'Cosmos'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range("D24") '11.06'
$ws.Range("E24").Value = '  -3.85%  '

# Row 25
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range("D25") '2.236'
$ws.Range("E25").Value = '  -2.02%  '

# Row 26
Set-TextValue $ws.Range("D26") '2.547'
$ws.Range("E26").Value = '  +2.12%  '

# Row 27
Set-TextValue $ws.Range("D27") '2.068.79'
$ws.Range("E27").Value = '  -0.52%  '

# Row 28
Set-TextValue $ws.Range("D28") '157.97'
$ws.Range("E28").Value = '  -0.09%  '

# Row 29
Set-TextValue $ws.Range("D29") '20.42'
$ws.Range("E29").Value = '  -1.33%  '

# Row 30
Set-TextValue $ws.Range("D30") '124.34'
$ws.Range("E30").Value = '  -1.28%  '

# Row 31
Set-TextValue $ws.Range("D31") '0.1055'
$ws.Range("E31").Value = '  -0.30%  '

# Row 32
Set-TextValue $ws.Range("D32") '1.037'
$ws.Range("E32").Value = '  -1.80%  '

# Row 33
Set-TextValue $ws.Range("D33") '5.611'
$ws.Range("E33").Value = '  +0.67%  '

# Row 34
Set-TextValue $ws.Range("D34") '3.598'
$ws.Range("E34").Value = '  +0.44%  '

# Row 35
Set-TextValue $ws.Range("D35") '9.626'
$ws.Range("E35").Value = '  +2.56%  '

# Row 36
Set-TextValue $ws.Range("D36") '0.06539'
$ws.Range("E36").Value = '  -0.20%  '

# Row 37
Set-TextValue $ws.Range("D37") '0.02407'
$ws.Range("E37").Value = '  -0.13%  '

# Row 38
Set-TextValue $ws.Range("D38") '0.2169'
$ws.Range("E38").Value = '  -0.90%  '

# Row 39
Set-TextValue $ws.Range("D39") '1.203'
$ws.Range("E39").Value = '  +0.64%  '

# Row 40
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range("D40") '0.6397'
$ws.Range("E40").Value = '  +0.53%  '

# Row 41
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range("D41") '1.235'
$ws.Range("E41").Value = '  -4.21%  '

# Row 42
$ws.Range("B42").Value = 'InternetComputer(DFINITY)'
$ws.Range("C42").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range("D42") '4.854'
$ws.Range("E42").Value = '  -0.77%  '

# Row 43
$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range("D43") '11.14'
$ws.Range("E43").Value = '  -3.77%  '

# Row 44
Set-TextValue $ws.Range("D44") '0.6065'
$ws.Range("E44").Value = '  +1.01%  '

# Row 45
Set-TextValue $ws.Range("D45") '13.06'
$ws.Range("E45").Value = '  -0.55%  '

# Row 46
$ws.Range("E46").Value = '  +0.05%  '

# Row 47
Set-TextValue $ws.Range("D47") '3.655'
$ws.Range("E47").Value = '  -0.57%  '

# Row 48
Set-TextValue $ws.Range("D48") '1.989'
$ws.Range("E48").Value = '  -0.23%  '

# Row 49
Set-TextValue $ws.Range("D49") '1.205'
$ws.Range("E49").Value = '  -1.16%  '

# Row 50
Set-TextValue $ws.Range("D50") '119.99'
$ws.Range("E50").Value = '  -0.85%  '

# Row 51
Set-TextValue $ws.Range("D51") '79.06'
$ws.Range("E51").Value = '  +0.40%  '
